# Weekly update: a new price record for "Poroto verde" (Macroferia Regional
# de Talca) is added as the newest row. The sheet's rows are stored most-
# recent-first starting at row 37 (rows 2-36 hold a separate, already up to
# date block), so the new record is inserted at row 37, pushing the former
# rows 37-150 down to 38-151 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 37; everything below (old rows 37-150) shifts
# down to 38-151, carrying its formatting (incl. the date style on column D).
$ws.Rows("37:37").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A37").Value = 5
$ws.Range("B37").Value = "Macroferia Regional de Talca"
$ws.Range("C37").Value = "Maule"
$ws.Range("D37").Value = 44622
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = 100112031
$ws.Range("G37").Value = "Poroto verde"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 80
$ws.Range("K37").Value = 30000
$ws.Range("L37").Value = 30000
$ws.Range("M37").Value = 30000
$ws.Range("N37").Value = "`$/saco 25 kilos"
$ws.Range("O37").Value = "Región del Maule"
$ws.Range("P37").Value = 1200
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
